$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Esplanade, Lakes Entrance VIC 3909, Australia'
$ws.Range("B2").Value = -37.8806917
$ws.Range("C2").Value = 147.987898
$ws.Range("D2").Value = 'East Gippsland (S)'

$ws.Range("A3").Value = 'North Arm Tourist Park, Marine Parade, Lakes Entrance VIC 3909, Australia'
$ws.Range("B3").Value = -37.8808079
$ws.Range("C3").Value = 147.9766688
$ws.Range("D3").Value = 'East Gippsland (S)'

$ws.Range("A4").Value = 'The Esplanade Resort & Spa, 1 Esplanade, Lakes Entrance VIC 3909, Australia'
$ws.Range("B4").Value = -37.8826838
$ws.Range("C4").Value = 147.9746235
$ws.Range("D4").Value = 'East Gippsland (S)'

$ws.Range("A5").Value = 'Wyanga Park Winery, Blairs Road, Lakes Entrance VIC 3909, Australia'
$ws.Range("B5").Value = -37.8552364
$ws.Range("C5").Value = 147.9848786
$ws.Range("D5").Value = 'East Gippsland (S)'

$ws.Range("A6").Value = 'European Bier Cafe, 120 Exhibition Street, Melbourne VIC 3000, Australia'
$ws.Range("B6").Value = -37.8127778
$ws.Range("C6").Value = 144.9712083
$ws.Range("D6").Value = 'Melbourne (C)'

$ws.Range("A7").Value = 'Sackville House, 27 Flinders Lane, Melbourne VIC 3000, Australia'
$ws.Range("B7").Value = -37.814891
$ws.Range("C7").Value = 144.9732909
$ws.Range("D7").Value = 'Melbourne (C)'

$ws.Range("A8").Value = 'Central Avenue, Moorabbin VIC 3189, Australia'
$ws.Range("B8").Value = -37.9349928
$ws.Range("C8").Value = 145.0371527
$ws.Range("D8").Value = 'Kingston (C) (Vic.)'

$ws.Range("A9").Value = 'Q3, 1 Southbank Boulevard, Southbank VIC 3006, Australia'
$ws.Range("B9").Value = -37.8209645
$ws.Range("C9").Value = 144.9625172
$ws.Range("D9").Value = 'Melbourne (C)'

$ws.Range("A10").Value = 'Crown Entertainment Complex, 8 Whiteman Street, Southbank VIC 3006, Australia'
$ws.Range("B10").Value = -37.8218482
$ws.Range("C10").Value = 144.96027
$ws.Range("D10").Value = 'Melbourne (C)'

$ws.Range("A11").Value = 'North Arm Tourist Park, Marine Parade, Lakes Entrance VIC 3909, Australia'
$ws.Range("B11").Value = -37.8803903
$ws.Range("C11").Value = 147.9783877
$ws.Range("D11").Value = 'East Gippsland (S)'

$ws.Range("A12").Value = 'Cunningham Arm, Elwyn Street, Kalimna VIC 3909, Australia'
$ws.Range("B12").Value = -37.8821077
$ws.Range("C12").Value = 147.9946854
$ws.Range("D12").Value = 'East Gippsland (S)'

$ws.Range("A13").Value = 'KFC, Esplanade, Lakes Entrance VIC 3909, Australia'
$ws.Range("B13").Value = -37.8812255
$ws.Range("C13").Value = 147.9848975
$ws.Range("D13").Value = 'East Gippsland (S)'

$ws.Range("A14").Value = 'Woolworths, Church Street, Lakes Entrance VIC 3909, Australia'
$ws.Range("B14").Value = -37.8797527
$ws.Range("C14").Value = 147.9870802
$ws.Range("D14").Value = 'East Gippsland (S)'

$ws.Range("A15").Value = 'Ritchie''s IGA, Hamilton Place, Mount Waverley VIC 3149, Australia'
$ws.Range("B15").Value = -37.8763889
$ws.Range("C15").Value = 145.1286111
$ws.Range("D15").Value = 'Monash (C)'

$ws.Range("A16").Value = 'The Links Shopping Centre, Luntar Road, Oakleigh South VIC 3167, Australia'
$ws.Range("B16").Value = -37.9219217
$ws.Range("C16").Value = 145.084782
$ws.Range("D16").Value = 'Monash (C)'

$ws.Range("A17").Value = 'Spink Street, Brighton VIC 3186, Australia'
$ws.Range("B17").Value = -37.8978842
$ws.Range("C17").Value = 145.0043056
$ws.Range("D17").Value = 'Bayside (C)'

$ws.Range("A18").Value = 'Aldi, Bay Road, Sandringham VIC 3190, Australia'
$ws.Range("B18").Value = -37.9558421
$ws.Range("C18").Value = 145.0334377
$ws.Range("D18").Value = 'Bayside (C)'

$ws.Range("A19").Value = 'Esplanade, Lakes Entrance VIC 3909, Australia'
$ws.Range("B19").Value = -37.8815345
$ws.Range("C19").Value = 147.9803919
$ws.Range("D19").Value = 'East Gippsland (S)'

$ws.Range("A20").Value = 'The Esplanade Resort & Spa, 1 Esplanade, Lakes Entrance VIC 3909, Australia'
$ws.Range("B20").Value = -37.8826838
$ws.Range("C20").Value = 147.9746235
$ws.Range("D20").Value = 'East Gippsland (S)'

$ws.Range("A21").Value = 'The Esplanade Resort & Spa, 1 Esplanade, Lakes Entrance VIC 3909, Australia'
$ws.Range("B21").Value = -37.8826838
$ws.Range("C21").Value = 147.9746235
$ws.Range("D21").Value = 'East Gippsland (S)'

$ws.Range("A22").Value = 'Bunnings Warehouse, 23-27 Nepean Highway, Mentone VIC 3194, Australia'
$ws.Range("B22").Value = -37.9768067
$ws.Range("C22").Value = 145.0694988
$ws.Range("D22").Value = 'Kingston (C) (Vic.)'

$ws.Range("A23").Value = 'Bunnings Warehouse, 23-27 Nepean Highway, Mentone VIC 3194, Australia'
$ws.Range("B23").Value = -37.9768067
$ws.Range("C23").Value = 145.0694988
$ws.Range("D23").Value = 'Kingston (C) (Vic.)'
